# Refresh the cryptos table with the latest scraped price/volume snapshot,
# matching the upstream "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe marks the literal as text (the same trick used when
# typing directly into Excel), so numeric-looking values such as "1.00" or
# "0.999" are written verbatim instead of being normalized into a number and
# losing their formatting (trailing zeros, etc).

$ws.Range("D2").Value = "'68.183.45"
$ws.Range("E2").Value = "'  +0.79%  "

$ws.Range("D3").Value = "'3.680.45"
$ws.Range("E3").Value = "'  +0.00%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.11%  "

$ws.Range("D5").Value = "'602.59"
$ws.Range("E5").Value = "'  +5.57%  "

$ws.Range("D6").Value = "'192.73"
$ws.Range("E6").Value = "'  +8.54%  "

$ws.Range("D7").Value = "'0.624"
$ws.Range("E7").Value = "'  +0.38%  "

$ws.Range("E8").Value = "'  -0.21%  "

$ws.Range("D9").Value = "'0.709"
$ws.Range("E9").Value = "'  +0.92%  "

$ws.Range("D10").Value = "'58.32"
$ws.Range("E10").Value = "'  +13.25%  "

$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "'  -3.89%  "

$ws.Range("D12").Value = "'0.0000277"
$ws.Range("E12").Value = "'  -3.10%  "

$ws.Range("D13").Value = "'10.29"
$ws.Range("E13").Value = "'  -0.65%  "

$ws.Range("D14").Value = "'4.258.41"
$ws.Range("E14").Value = "'  -0.38%  "

$ws.Range("D15").Value = "'3.672.45"
$ws.Range("E15").Value = "'  -0.35%  "

$ws.Range("E16").Value = "'  +1.11%  "

$ws.Range("D17").Value = "'19.06"
$ws.Range("E17").Value = "'  -0.78%  "

$ws.Range("D18").Value = "'1.13"
$ws.Range("E18").Value = "'  +1.51%  "

$ws.Range("D19").Value = "'67.965.75"
$ws.Range("E19").Value = "'  +0.73%  "

$ws.Range("D20").Value = "'12.57"
$ws.Range("E20").Value = "'  -1.12%  "

$ws.Range("D21").Value = "'403.56"
$ws.Range("E21").Value = "'  +0.09%  "

$ws.Range("D22").Value = "'4.47"
$ws.Range("E22").Value = "'  +1.53%  "

$ws.Range("D23").Value = "'88.46"
$ws.Range("E23").Value = "'  +1.00%  "

$ws.Range("D24").Value = "'11.38"
$ws.Range("E24").Value = "'  +6.25%  "

$ws.Range("D25").Value = "'2.98"
$ws.Range("E25").Value = "'  -1.16%  "

$ws.Range("D26").Value = "'12.66"
$ws.Range("E26").Value = "'  +0.38%  "

$ws.Range("D27").Value = "'6.03"
$ws.Range("E27").Value = "'  -0.12%  "

$ws.Range("D28").Value = "'3.70"
$ws.Range("E28").Value = "'  -1.72%  "

$ws.Range("D29").Value = "'9.39"
$ws.Range("E29").Value = "'  -0.38%  "

$ws.Range("D30").Value = "'32.14"
$ws.Range("E30").Value = "'  -0.41%  "

$ws.Range("D31").Value = "'7.60"
$ws.Range("E31").Value = "'  +2.85%  "

$ws.Range("B32").Value = "'OKB"
$ws.Range("C32").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "'68.37"
$ws.Range("E32").Value = "'  +5.99%  "

$ws.Range("B33").Value = "'Cosmos"
$ws.Range("C33").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'12.46"
$ws.Range("E33").Value = "'  +0.99%  "

$ws.Range("B34").Value = "'InjectiveProtocol"
$ws.Range("C34").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'45.75"
$ws.Range("E34").Value = "'  +6.73%  "

$ws.Range("E35").Value = "'  +2.76%  "

$ws.Range("D36").Value = "'613.39"
$ws.Range("E36").Value = "'  +0.64%  "

$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "'  -0.05%  "

$ws.Range("D38").Value = "'0.400"
$ws.Range("E38").Value = "'  +1.93%  "

$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "'  -0.13%  "

$ws.Range("D40").Value = "'0.0₃0775"
$ws.Range("E40").Value = "'  -11.44%  "

$ws.Range("D41").Value = "'0.136"
$ws.Range("E41").Value = "'  +1.33%  "

$ws.Range("D42").Value = "'2.94"
$ws.Range("E42").Value = "'  -0.96%  "

$ws.Range("D43").Value = "'0.0430"
$ws.Range("E43").Value = "'  +0.02%  "

$ws.Range("D44").Value = "'2.57"
$ws.Range("E44").Value = "'  -6.15%  "

$ws.Range("D45").Value = "'2.832.96"
$ws.Range("E45").Value = "'  +1.56%  "

$ws.Range("E46").Value = "'  +3.30%  "

$ws.Range("B47").Value = "'THORChain"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "'9.04"
$ws.Range("E47").Value = "'  -0.98%  "

$ws.Range("B48").Value = "'ApeXProtocol"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.19"
$ws.Range("E48").Value = "'  +4.20%  "

$ws.Range("D49").Value = "'144.69"
$ws.Range("E49").Value = "'  +4.68%  "

$ws.Range("E50").Value = "'  -0.93%  "

$ws.Range("E51").Value = "'  -11.79%  "
